$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 995.5
$ws.Range("I2").Value = 439.44446
$ws.Range("K2").Value = 439.44446
$ws.Range("M2").Value = -326.44446
$ws.Range("I9").Value = 8277.182000000001
$ws.Range("J9").Value = 94.75
$ws.Range("K9").Value = 8277.182000000001
$ws.Range("L9").Value = 94.75
$ws.Range("M9").Value = -8108.182000000001
$ws.Range("N9").Value = -432.75
$ws.Range("H40").Value = 1119.6072
$ws.Range("J40").Value = 1333.2222
$ws.Range("L40").Value = 1333.2222
$ws.Range("N40").Value = -1683.2222
$ws.Range("H53").Value = 264.3846
$ws.Range("I53").Value = 171.5
$ws.Range("J53").Value = 413
$ws.Range("K53").Value = 171.5
$ws.Range("L53").Value = 413
$ws.Range("M53").Value = 465.5
$ws.Range("N53").Value = -1687
$ws.Range("H92").Value = 5269.857
$ws.Range("I92").Value = 5778
$ws.Range("K92").Value = 5778
$ws.Range("M92").Value = -4530
$ws.Range("H107").Value = 7534.3
$ws.Range("J107").Value = 12899.75
$ws.Range("L107").Value = 12899.75
$ws.Range("N107").Value = -16739.75
$ws.Range("H135").Value = 3426.3333
$ws.Range("I135").Value = 3311.6
$ws.Range("K135").Value = 29804.4
$ws.Range("M135").Value = -27269.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2307.5
$ws.Range("I2").Value = 2424.077
$ws.Range("J2").Value = 2139.111
$ws.Range("K2").Value = 2424.077
$ws.Range("L2").Value = 2139.111
$ws.Range("M2").Value = -2311.077
$ws.Range("N2").Value = -2365.111
$ws.Range("H32").Value = 2358.9
$ws.Range("I32").Value = 1956.5352
$ws.Range("K32").Value = 1956.5352
$ws.Range("M32").Value = -1669.5352
$ws.Range("H45").Value = 23251.584
$ws.Range("J45").Value = 22001.6
$ws.Range("L45").Value = 22001.6
$ws.Range("N45").Value = -22755.6
$ws.Range("H61").Value = 3098.889
$ws.Range("I61").Value = 2682.0833
$ws.Range("K61").Value = 2682.0833
$ws.Range("M61").Value = -2470.0833
$ws.Range("H101").Value = 23533.666
$ws.Range("J101").Value = 23533.666
$ws.Range("L101").Value = 23533.666
$ws.Range("N101").Value = -30023.666
$ws.Range("H106").Value = 197800
$ws.Range("J106").Value = 197800
$ws.Range("L106").Value = 197800
$ws.Range("N106").Value = -200324
$ws.Range("H116").Value = 2307.5
$ws.Range("I116").Value = 2424.077
$ws.Range("J116").Value = 2139.111
$ws.Range("K116").Value = 2424.077
$ws.Range("L116").Value = 2139.111
$ws.Range("M116").Value = -130.0770000000002
$ws.Range("N116").Value = -6727.111
$ws.Range("H132").Value = 2728.8333
$ws.Range("I132").Value = 2456.9524
$ws.Range("J132").Value = 4632
$ws.Range("K132").Value = 7370.8572
$ws.Range("L132").Value = 13896
$ws.Range("M132").Value = -4840.8572
$ws.Range("N132").Value = -18956
$ws.Range("H136").Value = 3098.889
$ws.Range("I136").Value = 2682.0833
$ws.Range("K136").Value = 8046.249899999999
$ws.Range("M136").Value = -5496.249899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2307.5
$ws.Range("I3").Value = 2424.077
$ws.Range("J3").Value = 2139.111
$ws.Range("K3").Value = 2424.077
$ws.Range("L3").Value = 2139.111
$ws.Range("M3").Value = -2310.077
$ws.Range("N3").Value = -2367.111
$ws.Range("H92").Value = 60500
$ws.Range("J92").Value = 60500
$ws.Range("L92").Value = 60500
$ws.Range("N92").Value = -65492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 738.44446
$ws.Range("I16").Value = 768.25
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 768.25
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -481.25
$ws.Range("N16").Value = -1074
$ws.Range("H94").Value = 1233.9445
$ws.Range("I94").Value = 1307.0769
$ws.Range("K94").Value = 1307.0769
$ws.Range("M94").Value = -856.0769
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H113").Value = 738.44446
$ws.Range("I113").Value = 768.25
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 768.25
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = 1401.75
$ws.Range("N113").Value = -4840
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 196.5
$ws.Range("J7").Value = 164
$ws.Range("L7").Value = 492
$ws.Range("N7").Value = -716
$ws.Range("H55").Value = 1691.4
$ws.Range("J55").Value = 1691.4
$ws.Range("L55").Value = 5074.200000000001
$ws.Range("N55").Value = -5428.200000000001
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H75").Value = 8163.3335
$ws.Range("J75").Value = 8163.3335
$ws.Range("L75").Value = 24490.0005
$ws.Range("N75").Value = -26486.0005
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H78").Value = 8163.3335
$ws.Range("J78").Value = 8163.3335
$ws.Range("L78").Value = 73470.0015
$ws.Range("N78").Value = -83454.0015
$ws.Range("H107").Value = 1071.4286
$ws.Range("J107").Value = 1452
$ws.Range("L107").Value = 4356
$ws.Range("N107").Value = -8196
$ws.Range("H141").Value = 173402
$ws.Range("I141").Value = 8075.8
$ws.Range("K141").Value = 24227.4
$ws.Range("M141").Value = -19047.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9948.875
$ws.Range("I57").Value = 4898.75
$ws.Range("K57").Value = 4898.75
$ws.Range("M57").Value = -4078.75
$ws.Range("H102").Value = 1371.9714
$ws.Range("I102").Value = 758.2258
$ws.Range("J102").Value = 6128.5
$ws.Range("K102").Value = 758.2258
$ws.Range("L102").Value = 6128.5
$ws.Range("M102").Value = 863.7742
$ws.Range("N102").Value = -9372.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18834
$ws.Range("I7").Value = 20800.8
$ws.Range("K7").Value = 20800.8
$ws.Range("M7").Value = -20688.8
$ws.Range("H22").Value = 919.9
$ws.Range("I22").Value = 972.1818
$ws.Range("J22").Value = 856
$ws.Range("K22").Value = 972.1818
$ws.Range("L22").Value = 856
$ws.Range("M22").Value = -677.1818
$ws.Range("N22").Value = -1446
$ws.Range("H27").Value = 919.9
$ws.Range("I27").Value = 972.1818
$ws.Range("J27").Value = 856
$ws.Range("K27").Value = 972.1818
$ws.Range("L27").Value = 856
$ws.Range("M27").Value = -865.1818
$ws.Range("N27").Value = -1070
$ws.Range("H46").Value = 3007.5
$ws.Range("I46").Value = 2620.875
$ws.Range("J46").Value = 3523
$ws.Range("K46").Value = 2620.875
$ws.Range("L46").Value = 3523
$ws.Range("M46").Value = -2432.875
$ws.Range("N46").Value = -3899
$ws.Range("H126").Value = 18834
$ws.Range("I126").Value = 20800.8
$ws.Range("K126").Value = 62402.39999999999
$ws.Range("M126").Value = -59932.39999999999
$ws.Range("H136").Value = 4014.2964
$ws.Range("I136").Value = 3791.2273
$ws.Range("K136").Value = 11373.6819
$ws.Range("M136").Value = -8823.6819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 335133.34
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 2700
$ws.Range("K3").Value = 1000000
$ws.Range("L3").Value = 2700
$ws.Range("M3").Value = -999886
$ws.Range("N3").Value = -2928
$ws.Range("H9").Value = 49999
$ws.Range("K9").Value = 49999
$ws.Range("M9").Value = -49859
$ws.Range("H100").Value = 1200.2572
$ws.Range("I100").Value = 1049.381
$ws.Range("K100").Value = 2098.762
$ws.Range("M100").Value = -1557.762
$ws.Range("H136").Value = 2209.4595
$ws.Range("I136").Value = 1919.1724
$ws.Range("J136").Value = 3261.75
$ws.Range("K136").Value = 5757.5172
$ws.Range("L136").Value = 9785.25
$ws.Range("M136").Value = -3207.5172
$ws.Range("N136").Value = -14885.25
